$d = $word.ActiveDocument

$replacements = @(
    @{old = "405÷8=50, 5"; new = "990÷5=198, 0"},
    @{old = "852÷3=284, 0"; new = "517÷8=64, 5"},
    @{old = "146÷9=16, 2"; new = "922÷7=131, 5"},
    @{old = "796÷6=132, 4"; new = "394÷8=49, 2"},
    @{old = "321÷2=160, 1"; new = "383÷8=47, 7"},
    @{old = "358÷4=89, 2"; new = "354÷9=39, 3"},
    @{old = "407÷6=67, 5"; new = "287÷9=31, 8"},
    @{old = "760÷4=190, 0"; new = "749÷2=374, 1"},
    @{old = "364÷9=40, 4"; new = "636÷7=90, 6"},
    @{old = "564÷5=112, 4"; new = "319÷9=35, 4"},
    @{old = "724÷7=103, 3"; new = "181÷2=90, 1"},
    @{old = "463÷4=115, 3"; new = "965÷5=193, 0"},
    @{old = "351÷5=70, 1"; new = "648÷5=129, 3"},
    @{old = "467÷7=66, 5"; new = "938÷7=134, 0"},
    @{old = "250÷2=125, 0"; new = "428÷6=71, 2"},
    @{old = "637÷2=318, 1"; new = "157÷7=22, 3"},
    @{old = "560÷2=280, 0"; new = "340÷7=48, 4"},
    @{old = "372÷7=53, 1"; new = "193÷5=38, 3"},
    @{old = "209÷3=69, 2"; new = "314÷6=52, 2"},
    @{old = "389÷6=64, 5"; new = "452÷2=226, 0"},
    @{old = "797÷8=99, 5"; new = "537÷6=89, 3"},
    @{old = "371÷3=123, 2"; new = "777÷9=86, 3"},
    @{old = "178÷5=35, 3"; new = "550÷2=275, 0"},
    @{old = "730÷7=104, 2"; new = "354÷7=50, 4"},
    @{old = "114÷4=28, 2"; new = "265÷3=88, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
